$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 121; existing rows 121-178 shift down to 122-179.
$ws.Rows("121:121").Insert()

# Populate the newly inserted row 121 with the new weekly record.
$ws.Range("A121").Value = 4
$ws.Range("B121").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C121").Value = "Los Lagos"
$ws.Range("D121").Value = 44523
$ws.Range("E121").Value = 10
$ws.Range("F121").Value = 100112003
$ws.Range("G121").Value = "Ajo"
$ws.Range("H121").Value = "Chino"
$ws.Range("I121").Value = "Primera"
$ws.Range("J121").Value = 240
$ws.Range("K121").Value = 21000
$ws.Range("L121").Value = 21000
$ws.Range("M121").Value = 21000
$ws.Range("N121").Value = "$/caja 10 kilos"
$ws.Range("O121").Value = "China"
$ws.Range("P121").Value = 2100
$ws.Range("Q121").Value = 10
$ws.Range("R121").Value = "Hortaliza"
